# Adds the "Most likely unigrams / bigrams / trigrams" table (F35:K41)
# plus a handful of empty, styled cells below it (F45:F48) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (35) ------------------------------------------------
$ws.Range("F35").Value = "Most likely unigrams"
$ws.Range("H35").Value = "Most likely bigrams"
$ws.Range("J35").Value = "Most likely trigrams"

# --- Column-header row (36) ------------------------------------------
$ws.Range("F36").Value = "unigram"
$ws.Range("G36").Value = "FPMI"
$ws.Range("H36").Value = "bigrams"
$ws.Range("I36").Value = "FPMI"
$ws.Range("J36").Value = "trigram"
$ws.Range("K36").Value = "FPMI"

# --- Data rows (37-41) -------------------------------------------------
$unigrams  = @("the", ".", "and", "of", "to")
$uniCounts = @(50876, 49836, 32774, 29508, 22812)

$bigrams   = @("of the ", "in the", ". the ", ". i ", ". and ")
$biCounts  = @(6651, 4099, 2953, 2765, 2698)

$trigrams  = @(". and the", ". it is ", "of the lord ", "the son of ", "out of the ")
$triCounts = @(365, 329, 315, 272, 261)

for ($i = 0; $i -lt 5; $i++) {
    $row = 37 + $i

    $ws.Cells.Item($row, 6).Value = $unigrams[$i]
    $ws.Cells.Item($row, 7).Value = $uniCounts[$i]

    $ws.Cells.Item($row, 8).Value = $bigrams[$i]
    $ws.Cells.Item($row, 9).Value = $biCounts[$i]

    $ws.Cells.Item($row, 10).Value = $trigrams[$i]
    $ws.Cells.Item($row, 11).Value = $triCounts[$i]
}

# --- Styled (font size 10, Helvetica) cells --------------------------
$styledRange = $ws.Range("F37:F41,H37:H41,J37:J41,F45:F48")
$styledRange.Font.Name = "Helvetica"
$styledRange.Font.Size = 10

# --- Column width for J -----------------------------------------------
$ws.Columns.Item(10).ColumnWidth = 17

# --- View / selection state -------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("K41").Select()
